$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 30: codice_particella=2727/1, codice_comune_catastale=189 ----
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "2727/1"
$ws.Cells.Item(30, 3).Value = 189

# ---- Row 31: codice_particella=.256, codice_comune_catastale=231 ----
$ws.Cells.Item(31, 1).Value = 29

# ".256" parses as a numeric literal, so write it through a temporary
# text-formatted helper cell and paste back only the value. This keeps
# the content as literal text without leaving a stray number-format
# style applied to the destination cell.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = ".256"
$helper.Copy()
$ws.Range("B31").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item(31, 3).Value = 231

# ---- Match column-A formatting (border + centered/bold style) used by
# the rest of the table for the two newly added rows ----
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
